$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.116.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.72%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.649.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.16%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.49%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.68%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.644.07"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.18%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.627"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.71%  "

# Row 9
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.710"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.54%  "

# Row 11
$ws.Range("E11").Value = "  -8.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000290"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.77%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.03%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.224.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.18%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.647.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.40%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -8.96%  "

# Row 18
$ws.Range("E18").Value = "  -2.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.89%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.870.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.82%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "408.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.79%  "

# Row 23
$ws.Range("E23").Value = "  -4.49%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.60%  "

# Row 25
$ws.Range("E25").Value = "  -9.11%  "

# Row 26
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.18%  "

# Row 27
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.81%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.27%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.31%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.65%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.10%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -14.28%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.79%  "

# Row 34
$ws.Range("E34").Value = "  -6.83%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "64.36"
$ws.Range("D35").Style = "Normal"

# Row 36
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "42.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -11.62%  "

# Row 37
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "592.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.52%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0876"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.19%  "

# Row 39
$ws.Range("E39").Value = "  +0.02%  "

# Row 40
$ws.Range("E40").Value = "  -8.73%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.10%  "

# Row 42
$ws.Range("E42").Value = "  -6.23%  "

# Row 43
$ws.Range("E43").Value = "  -7.16%  "

# Row 44
$ws.Range("E44").Value = "  -7.69%  "

# Row 45
$ws.Range("E45").Value = "  -7.02%  "

# Row 46
$ws.Range("E46").Value = "  -11.68%  "

# Row 47
$ws.Range("E47").Value = "  -3.92%  "

# Row 48
$ws.Range("E48").Value = "  -6.62%  "

# Row 49
$ws.Range("E49").Value = "  -10.53%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.710.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.75%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.81%  "
